$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.067.34"
$ws.Range("E2").Value = "  -5.35%  "
$ws.Range("D3").Value = "3.227.11"
$ws.Range("E3").Value = "  -8.60%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'582.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.19%  "
$ws.Range("D6").Value = "'151.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -12.43%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "3.220.76"
$ws.Range("E8").Value = "  -8.64%  "
$ws.Range("D9").Value = "'0.542"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -11.40%  "
$ws.Range("E10").Value = "  -13.15%  "
$ws.Range("D11").Value = "'6.72"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.94%  "
$ws.Range("D12").Value = "'0.501"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -15.07%  "
$ws.Range("D13").Value = "'38.09"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -18.34%  "
$ws.Range("E14").Value = "  -12.00%  "
$ws.Range("D15").Value = "3.743.27"
$ws.Range("E15").Value = "  -8.69%  "
$ws.Range("D16").Value = "66.963.47"
$ws.Range("E16").Value = "  -5.55%  "
$ws.Range("D17").Value = "3.226.49"
$ws.Range("E17").Value = "  -8.64%  "
$ws.Range("D18").Value = "'538.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -12.15%  "
$ws.Range("E19").Value = "  -5.96%  "
$ws.Range("E20").Value = "  -15.63%  "
$ws.Range("D21").Value = "'15.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -15.24%  "
$ws.Range("D22").Value = "'0.757"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -14.72%  "
$ws.Range("D23").Value = "'7.74"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -14.02%  "
$ws.Range("D24").Value = "'85.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -12.98%  "
$ws.Range("D25").Value = "'13.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -14.83%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "'3.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -16.67%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'29.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -13.57%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'8.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -12.60%  "
$ws.Range("D30").Value = "'2.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -18.09%  "
$ws.Range("D31").Value = "'2.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -14.40%  "
$ws.Range("E32").Value = "  -13.35%  "
$ws.Range("D33").Value = "'542.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -10.87%  "
$ws.Range("E34").Value = "  -20.17%  "
$ws.Range("D35").Value = "'5.68"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -17.37%  "
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("D37").Value = "'53.53"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.12%  "
$ws.Range("D38").Value = "'0.0434"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.57%  "
$ws.Range("E39").Value = "  -15.87%  "
$ws.Range("D40").Value = "'0.0841"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -16.36%  "
$ws.Range("D41").Value = "'0.126"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -13.25%  "
$ws.Range("D42").Value = "2.927.58"
$ws.Range("E42").Value = "  -13.30%  "
$ws.Range("D43").Value = "'2.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -27.52%  "
$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").Value = "0.0₃0584"
$ws.Range("E44").Value = "  -21.51%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "'0.259"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -17.30%  "
$ws.Range("D46").Value = "'2.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -20.10%  "
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("D48").Value = "'25.97"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -19.44%  "
$ws.Range("E49").Value = "  -18.80%  "
$ws.Range("E50").Value = "  -13.56%  "
$ws.Range("D51").Value = "'123.64"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.55%  "
